# Add team record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AD1 = Wins, AE1 = Losses, AF1 = Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (A1) so the new
# headers match the bold/centered/bordered look of the rest of row 1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Data rows 2-51: Wins = 90, Losses = 72, Ties = 0
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 90   # AD
    $ws.Cells.Item($r, 31).Value = 72   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
